$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "wind" and "solar" resource map columns with their color palettes
$ws.Range("E1").Value = "wind"
$ws.Range("F1").Value = "solar"

$ws.Range("F2").Value = "#c3eff1"
$ws.Range("F3").Value = "#90c2ec"
$ws.Range("F4").Value = "#4a9ae4"
$ws.Range("F5").Value = "#4b69c5"
$ws.Range("F6").Value = "#3b3aa5"

$ws.Range("E2").Value = "#feeed5"
$ws.Range("E3").Value = "#fdd3aa"
$ws.Range("E4").Value = "#ffbc85"
$ws.Range("E5").Value = "#f8a462"
$ws.Range("E6").Value = "#f19139"

# Update selection to match the new active cell
$ws.Range("E6").Select()
